$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.770.82'
$ws.Range('E2').Value = '  -0.77%  '

$ws.Range('D3').Value = '3.888.57'
$ws.Range('E3').Value = '  -1.98%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').Value = '''597.25'
$ws.Range('E5').Value = '  +1.87%  '

$ws.Range('D6').Value = '''167.34'
$ws.Range('E6').Value = '  +10.37%  '

$ws.Range('D7').Value = '''0.668'
$ws.Range('E7').Value = '  -0.80%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = '''0.762'
$ws.Range('E9').Value = '  +2.48%  '

$ws.Range('D10').Value = '''0.178'
$ws.Range('E10').Value = '  +7.19%  '

$ws.Range('D11').Value = '''53.90'
$ws.Range('E11').Value = '  +1.52%  '

$ws.Range('D12').Value = '''0.0000322'
$ws.Range('E12').Value = '  +1.92%  '

$ws.Range('D13').Value = '''11.34'
$ws.Range('E13').Value = '  +5.92%  '

$ws.Range('D14').Value = '4.519.50'
$ws.Range('E14').Value = '  -1.93%  '

$ws.Range('D15').Value = '3.900.90'
$ws.Range('E15').Value = '  -2.26%  '

$ws.Range('D16').Value = '''20.98'
$ws.Range('E16').Value = '  +3.12%  '

$ws.Range('D17').Value = '''13.90'
$ws.Range('E17').Value = '  -0.28%  '

$ws.Range('D18').Value = '''1.21'
$ws.Range('E18').Value = '  -5.37%  '

$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '''0.129'
$ws.Range('E19').Value = '  -1.79%  '

$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '71.724.44'
$ws.Range('E20').Value = '  -0.79%  '

$ws.Range('D21').Value = '''434.75'
$ws.Range('E21').Value = '  +2.03%  '

$ws.Range('D22').Value = '''4.73'
$ws.Range('E22').Value = '  +1.19%  '

$ws.Range('D23').Value = '''94.27'
$ws.Range('E23').Value = '  -1.04%  '

$ws.Range('D24').Value = '''3.29'
$ws.Range('E24').Value = '  -4.03%  '

$ws.Range('D25').Value = '''13.84'
$ws.Range('E25').Value = '  -2.12%  '

$ws.Range('D26').Value = '''4.16'
$ws.Range('E26').Value = '  -5.95%  '

$ws.Range('D27').Value = '''10.98'
$ws.Range('E27').Value = '  -1.83%  '

$ws.Range('D28').Value = '''5.93'
$ws.Range('E28').Value = '  +0.22%  '

$ws.Range('D29').Value = '''10.19'
$ws.Range('E29').Value = '  -3.82%  '

$ws.Range('D30').Value = '''35.12'
$ws.Range('E30').Value = '  -2.83%  '

$ws.Range('D31').Value = '''8.07'
$ws.Range('E31').Value = '  +4.24%  '

$ws.Range('D32').Value = '''52.14'
$ws.Range('E32').Value = '  +5.19%  '

$ws.Range('D33').Value = '''13.59'
$ws.Range('E33').Value = '  +1.62%  '

$ws.Range('E34').Value = '  -3.51%  '

$ws.Range('D35').Value = '0.0₃0980'
$ws.Range('E35').Value = '  +15.32%  '

$ws.Range('D36').Value = '''68.29'
$ws.Range('E36').Value = '  +0.17%  '

$ws.Range('D37').Value = '''621.76'
$ws.Range('E37').Value = '  -8.42%  '

$ws.Range('D38').Value = '''0.419'
$ws.Range('E38').Value = '  -3.48%  '

$ws.Range('E39').Value = '  +0.06%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.24%  '

$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').Value = '''3.30'
$ws.Range('E41').Value = '  +0.31%  '

$ws.Range('D42').Value = '''0.142'
$ws.Range('E42').Value = '  -2.01%  '

$ws.Range('D43').Value = '''3.20'
$ws.Range('E43').Value = '  +39.42%  '

$ws.Range('D44').Value = '''0.0470'
$ws.Range('E44').Value = '  -2.77%  '

$ws.Range('D45').Value = '''10.27'
$ws.Range('E45').Value = '  -5.99%  '

$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').Value = '''2.63'
$ws.Range('E46').Value = '  -4.01%  '

$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '''0.143'
$ws.Range('E47').Value = '  -2.45%  '

$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '''3.34'
$ws.Range('E48').Value = '  -1.26%  '

$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').Value = '''2.83'
$ws.Range('E49').Value = '  -15.09%  '

$ws.Range('D50').Value = '2.864.83'
$ws.Range('E50').Value = '  +3.49%  '

$ws.Range('D51').Value = '''0.000274'
$ws.Range('E51').Value = '  +1.79%  '
